# "+ arquivos sobre python"
#
# Adds two new rows of data to the "downloads/extras" table at the bottom
# of the sheet (rows 23 and 24), which were previously empty placeholder
# cells. The SUM formulas in B15 and B26 recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write row 24 before row 23 so the new shared strings are appended to
# sharedStrings.xml in the same order as in the target workbook:
#   19 -> "Exemplos_de_Caminho_de_Arquivos_em_Python"
#   20 -> "downloadArquivosNet"
$ws.Range("A24").Value = "Exemplos_de_Caminho_de_Arquivos_em_Python"
$ws.Range("B24").Value = 1

$ws.Range("A23").Value = "downloadArquivosNet"
$ws.Range("B23").Value = 3

# Reflect the new selection / scroll position recorded in the workbook.
$ws.Range("A30").Select()

$wb.Save()
